$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 12:03"

# Update province statistics
$ws.Range("B4").Value = 63416
$ws.Range("C4").Value = 38331
$ws.Range("D4").Value = 16619
$ws.Range("E4").Value = 8466
$ws.Range("B5").Value = 50924
$ws.Range("C5").Value = 22881
$ws.Range("D5").Value = 22698
$ws.Range("E5").Value = 5345
$ws.Range("B6").Value = 17520
$ws.Range("C6").Value = 7036
$ws.Range("D6").Value = 8637
$ws.Range("E6").Value = 1847
$ws.Range("B7").Value = 16144
$ws.Range("C7").Value = 5862
$ws.Range("D7").Value = 7635
$ws.Range("E7").Value = 2647
$ws.Range("B9").Value = 12236
$ws.Range("C9").Value = 7679
$ws.Range("D9").Value = 3276
$ws.Range("E9").Value = 1281
$ws.Range("B10").Value = 9097
$ws.Range("C10").Value = 6802
$ws.Range("D10").Value = 1713
$ws.Range("E10").Value = 582
$ws.Range("B13").Value = 5231
$ws.Range("C13").Value = 2790
$ws.Range("D13").Value = 1653
$ws.Range("E13").Value = 788
$ws.Range("B15").Value = 4966
$ws.Range("C15").Value = 2628
$ws.Range("D15").Value = 1862
$ws.Range("E15").Value = 476
$ws.Range("B16").Value = 3980
$ws.Range("C16").Value = 2396
$ws.Range("D16").Value = 1247
$ws.Range("E16").Value = 337
$ws.Range("B23").Value = 2865
$ws.Range("C23").Value = 2189
$ws.Range("D23").Value = 213
$ws.Range("E23").Value = 463
$ws.Range("B30").Value = 2310
$ws.Range("C30").Value = 953
$ws.Range("D30").Value = 1070
$ws.Range("E30").Value = 287
$ws.Range("C31").Value = 1223
$ws.Range("D31").Value = 865
$ws.Range("E31").Value = 143
$ws.Range("B33").Value = 2213
$ws.Range("C33").Value = 1696
$ws.Range("D33").Value = 318
$ws.Range("E33").Value = 199
$ws.Range("C59").Value = 104
$ws.Range("D59").Value = 13
